# Applies the cryptos list refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.293.09"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").Value = "2.545.19"
$ws.Range("E3").Value = "  -2.70%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "590.73"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "173.70"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.73%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("E8").Value = "  -0.33%  "

# Row 9
$ws.Range("D9").Value = "2.544.68"

# Row 10
$ws.Range("E10").Value = "  +0.85%  "

# Row 11
$ws.Range("E11").Value = "  +1.13%  "

# Row 12
$ws.Range("E12").Value = "  -0.74%  "

# Row 13
$ws.Range("E13").Value = "  -3.18%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.05"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.13%  "

# Row 15
$ws.Range("D15").Value = "3.011.23"
$ws.Range("E15").Value = "  -2.59%  "

# Row 16
$ws.Range("E16").Value = "  -0.68%  "

# Row 17
$ws.Range("D17").Value = "67.152.47"
$ws.Range("E17").Value = "  +0.06%  "

# Row 18
$ws.Range("D18").Value = "2.548.66"
$ws.Range("E18").Value = "  -2.65%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "8.05"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.29%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.43"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.93%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "356.47"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.22"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "

# Row 23
$ws.Range("E23").Value = "  +0.86%  "

# Row 24
$ws.Range("E24").Value = "  +4.99%  "

# Row 25
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "70.08"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.31%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.04"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -4.20%  "

# Row 28
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.17%  "

# Row 29
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.668.43"
$ws.Range("E29").Value = "  -2.85%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0998"
$ws.Range("E30").Value = "  +0.04%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "536.04"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.30%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.25"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +4.74%  "

# Row 33
$ws.Range("E33").Value = "  +0.67%  "

# Row 34
$ws.Range("E34").Value = "  -0.67%  "

# Row 35
$ws.Range("E35").Value = "  -0.88%  "

# Row 36
$ws.Range("E36").Value = "  +0.08%  "

# Row 37
$ws.Range("E37").Value = "  -0.43%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "157.05"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.39%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.82"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.55%  "

# Row 40
$ws.Range("E40").Value = "  +1.14%  "

# Row 41
$ws.Range("E41").Value = "  -2.01%  "

# Row 42
$ws.Range("E42").Value = "  +0.30%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.21"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.31%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.56"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +6.00%  "

# Row 45
$ws.Range("E45").Value = "  +0.03%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "39.90"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "151.13"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "

# Row 48
$ws.Range("E48").Value = "  -1.82%  "

# Row 49
$ws.Range("E49").Value = "  -4.99%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "3.72"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.18%  "

# Row 51
$ws.Range("E51").Value = "  +1.09%  "
